$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Add the new "national_id" column (M) ---
# Copy formatting (styles/borders) from column L into the new column M,
# matching rows 1-10 (the sheet's used range).
$ws.Range("L1:L10").Copy()
$ws.Range("M1:M10").PasteSpecial(-4122)

# Header for the new column
$ws.Range("M1").Value = "national_id"

# National id value for the first data row (row 2); row 3 is left blank
$ws.Range("M2").Value = "123"

# --- Column widths ---
# Re-apply column D:E width so it re-serializes through Excel's own
# rounding (23.4531 -> 23.5), matching the existing worksheet convention.
$ws.Range("D1:E1").ColumnWidth = 22.67

# Give the new column M the same width as the existing K:L columns.
$ws.Range("M1").ColumnWidth = $ws.Range("K1").ColumnWidth

"Done"
